# Updates odds values on the active worksheet to reflect the latest
# FlashScore data refresh (Atualizando o arquivo XLSX).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.91
$ws.Range("H2").Value = 3.2
$ws.Range("I2").Value = 4.75
$ws.Range("J2").Value = 2.75
$ws.Range("L2").Value = 5.5
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 6
$ws.Range("O2").Value = 1.57
$ws.Range("P2").Value = 2.25
$ws.Range("X2").Value = 7.5
$ws.Range("AG2").Value = 8.5
$ws.Range("AK2").Value = 41
$ws.Range("AW2").Value = 6

# Row 6
$ws.Range("O6").Value = 1.2
$ws.Range("P6").Value = 4.33

# Row 7
$ws.Range("G7").Value = 1.4
$ws.Range("N7").Value = 12
$ws.Range("S7").Value = 1.33
$ws.Range("T7").Value = 3.25
$ws.Range("U7").Value = 2
$ws.Range("V7").Value = 1.73
$ws.Range("W7").Value = 7
$ws.Range("Y7").Value = 9
$ws.Range("AC7").Value = 12
$ws.Range("AG7").Value = 17
$ws.Range("AI7").Value = 21
$ws.Range("AT7").Value = 3.25
$ws.Range("BA7").Value = 151

# Row 8
$ws.Range("G8").Value = 3.6
$ws.Range("I8").Value = 2.15
$ws.Range("Q8").Value = 2.1
$ws.Range("R8").Value = 1.7
$ws.Range("W8").Value = 10
$ws.Range("AJ8").Value = 19
$ws.Range("AO8").Value = 21
$ws.Range("AR8").Value = 101
$ws.Range("AU8").Value = 8.5

# Row 9
$ws.Range("G9").Value = 2.5
$ws.Range("L9").Value = 4
$ws.Range("Z9").Value = 23

# Row 12
$ws.Range("G12").Value = 3.9
$ws.Range("H12").Value = 3.3
$ws.Range("I12").Value = 1.8
$ws.Range("J12").Value = 4.75
$ws.Range("K12").Value = 2.1
$ws.Range("L12").Value = 2.6
$ws.Range("Q12").Value = 2.05
$ws.Range("R12").Value = 1.75
$ws.Range("X12").Value = 21
$ws.Range("Y12").Value = 15
$ws.Range("AC12").Value = 9
$ws.Range("AH12").Value = 8
$ws.Range("AJ12").Value = 15
$ws.Range("AX12").Value = 10
$ws.Range("AZ12").Value = 34
$ws.Range("BA12").Value = 51
$ws.Range("BB12").Value = 151
